$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" column (G) values: generator now derives K from a
# strike-run definition instead of the old raw Strike# count, and the
# per-row results were regenerated and re-written below (std/mean and
# s_vals are recalculated by the same pass, but do not change any other
# cell on this sheet).
$newK = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 2
    7 = 0
    8 = 2
    9 = 0
    10 = 3
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    24 = 2
    25 = 2
    26 = 0
    27 = 1
    28 = 2
    29 = 3
    30 = 1
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 2
    46 = 1
    47 = 0
    48 = 1
    49 = 1
    50 = 0
    51 = 2
    52 = 2
    53 = 4
    54 = 2
    55 = 1
    56 = 0
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 0
    63 = 2
    64 = 1
    65 = 1
    66 = 3
    67 = 3
    68 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

